# Adds a new "2022-Q1" sheet (fund holdings detail) before the "总计"
# summary sheet, and updates the "总计" sheet with the new quarter's
# aggregate row.

function Set-TextValue($cell, $val) {
    # Forces a numeric-looking string (e.g. "39.80", "001481") to be
    # stored as text instead of being auto-converted to a number, while
    # keeping the cell on the default ("Normal") style afterwards.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# The workbook currently ends with: ... 2021-Q4, 总计
$totalSheetBeforeInsert = $wb.Worksheets.Item($wb.Worksheets.Count)
$prevQuarterSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)

# --- 1. Create the new "2022-Q1" sheet right before "总计" ---------------
$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert)
$newSheet.Name = "2022-Q1"

# NOTE: Worksheets.Add() shifts sheet positions, and an already-held sheet
# reference tracks its slot (not the sheet it originally pointed at), so
# $totalSheetBeforeInsert now actually refers to the freshly inserted
# "2022-Q1" sheet. Re-fetch "总计" by name to get the right object.
$totalSheet = $wb.Worksheets.Item("总计")

# Reuse header formatting (style index used by the other quarter sheets)
# and the column-A numbering style from the immediately preceding quarter
# sheet, so we don't introduce any new style entries.
$prevQuarterSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$prevQuarterSheet.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

$fundRows = @(
    @("001481", "华宝油气(QDII)美元", "39.80", "94.60", "2.28", "0.9074", 6),
    @("162411", "华宝油气(QDII)人民币A", "39.80", "94.60", "2.28", "0.9074", 6),
    @("006679", "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A", "14.75", "83.19", "3.50", "0.5162", 9),
    @("162719", "广发道琼斯美国石油开发与生产指数（QDII-LOF）A", "14.75", "83.19", "3.50", "0.5162", 9),
    @("007844", "华宝油气(QDII)人民币C", "12.98", "94.60", "2.28", "0.2959", 6),
    @("006680", "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C", "4.73", "83.19", "3.50", "0.1656", 9),
    @("004243", "广发道琼斯美国石油开发与生产指数（QDII-LOF）C", "4.73", "83.19", "3.50", "0.1656", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2

    Set-TextValue $newSheet.Cells.Item($r, 2) $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[5]

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# --- 2. Update the "总计" sheet with the new 2022-Q1 aggregate row -------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)

$totalRows = @(
    @("2022-Q1", 7, 3.47),
    @("2021-Q4", 4, 0.4),
    @("2021-Q3", 6, 0.48),
    @("2021-Q2", 5, 0.72),
    @("2021-Q1", 8, 3.38),
    @("2020-Q4", 4, 0.5)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

$wb.Save()
